# 1. Justify the first paragraph (Format > Paragraph > Alignment: Justified).
$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)
$p1.Format.Alignment = 3

# 2. Fence off the boundary between the ". " run and the "It provides two
#    ways..." run with a temporary bookmark so that the text edits below
#    (which normally cause adjoining plain runs to be coalesced into one)
#    cannot merge across that boundary. The fence bookmark is removed
#    again once the new content is in place (removing a bookmark does not
#    trigger a re-merge).
$fenceRange = $d.Range(297, 297)
$d.Bookmarks.Add("FENCE1", $fenceRange)

# 3. Replace the run tail "GloVe, short for Global Vectors for Word
#    Representation" (which, in the original, was split across the
#    "_GoBack" bookmark) with an expanded discussion covering GloVe and
#    fastText, broken up into several runs exactly as the target revision
#    has it.
$rng = $d.Content
$rng.Find.Execute("GloVe, short for Global Vectors for Word Representation", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertPos = $rng.Start
$rng.Text = ""

$pieces = @(
    "Another popular model for obtaining word embeddings is ",
    "GloVe, short for Global Vectors for Word Representation",
    ", ",
    "which ",
    "is an extension of word2vec",
    ". It is a count-based model which ",
    "uses both global matrix factorization and local context window methods to obtain quality embeddings. fastText is yet another popular extension of word2vec. ",
    "This model represents each word as a bag of character n-grams. This allows capturing the meanings of shorter words and prefixes/suffixes. "
)

$cursor = $insertPos
foreach ($piece in $pieces) {
    $ip = $d.Range($cursor, $cursor)
    $ip.InsertAfter($piece)
    $cursor = $cursor + $piece.Length
}

$d.Bookmarks("FENCE1").Delete()

# 4. Re-insert the _GoBack bookmark immediately after the new text (i.e.
#    at the very end of the paragraph, before the paragraph mark). A
#    temporary placeholder character is used so the insertion point isn't
#    the very last position in the story when the bookmark is created
#    (adding a collapsed bookmark exactly at story-end otherwise expands
#    it to cover the whole range) -- the placeholder is removed right
#    after.
$dummyPos = $d.Range($cursor, $cursor)
$dummyPos.InsertAfter("Z")

$bmRange = $d.Range($cursor, $cursor)
$d.Bookmarks.Add("_GoBack", $bmRange)

$dummyRange = $d.Range($cursor, $cursor + 1)
$dummyRange.Text = ""
